# Auto-generated edit script: updates Leve profit calculations on multiple sheets
# (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR) to reflect refreshed market-board prices.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17: One for the Road
$ws.Range("H17").Value = 5556691.5
$ws.Range("J17").Value = 5556691.5
$ws.Range("L17").Value = 16670074.5
$ws.Range("N17").Value = -16670410.5

# Row 32: Automata for the People
$ws.Range("H32").Value = 2354.9
$ws.Range("I32").Value = 2056.6667
$ws.Range("J32").Value = 2482.7144
$ws.Range("K32").Value = 2056.6667
$ws.Range("L32").Value = 2482.7144
$ws.Range("M32").Value = -1730.6667
$ws.Range("N32").Value = -3134.7144

# Row 51: A Bile Business
$ws.Range("H51").Value = 5999
$ws.Range("I51").Value = 5999
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 5999
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -5515
$ws.Range("N51").ClearContents()

# Row 75: Tomes Roam on the Range
$ws.Range("H75").Value = 74666.336
$ws.Range("J75").Value = 74666.336
$ws.Range("L75").Value = 74666.336
$ws.Range("N75").Value = -76538.336

# Row 78: Field Trip to the Unknown (L)
$ws.Range("H78").Value = 74666.336
$ws.Range("J78").Value = 74666.336
$ws.Range("L78").Value = 223999.008
$ws.Range("N78").Value = -233359.008

# Row 111: An Eye for Healing
$ws.Range("H111").Value = 936.0909
$ws.Range("I111").Value = 874.5
$ws.Range("K111").Value = 2623.5
$ws.Range("M111").Value = 443.5

# Row 127: Liquid Competence
$ws.Range("H127").Value = 858.7273
$ws.Range("I127").Value = 794.6
$ws.Range("K127").Value = 2383.8
$ws.Range("M127").Value = 2576.2

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 1700
$ws.Range("I137").Value = 1523.8572
$ws.Range("J137").Value = 2624.75
$ws.Range("K137").Value = 4571.571599999999
$ws.Range("L137").Value = 7874.25
$ws.Range("M137").Value = -2021.571599999999
$ws.Range("N137").Value = -12974.25

# Row 138: All-night Crafting
$ws.Range("H138").Value = 3047.7144
$ws.Range("J138").Value = 3640.476
$ws.Range("L138").Value = 10921.428
$ws.Range("N138").Value = -21201.428

$ws = $wb.Worksheets.Item("ARM")
# Row 12: Strait Ain't the Gate
$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()

# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 8788.25
$ws.Range("I61").Value = 7480.706
$ws.Range("K61").Value = 7480.706
$ws.Range("M61").Value = -7268.706

# Row 122: Haste for High Durium
$ws.Range("H122").Value = 4900
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

# Row 136: Metal with Mettle
$ws.Range("H136").Value = 8788.25
$ws.Range("I136").Value = 7480.706
$ws.Range("K136").Value = 22442.118
$ws.Range("M136").Value = -19892.118

$ws = $wb.Worksheets.Item("BSM")
# Row 107: The Gold Experience
$ws.Range("H107").Value = 2575.963
$ws.Range("J107").Value = 5331.6665
$ws.Range("L107").Value = 5331.6665
$ws.Range("N107").Value = -9171.666499999999

# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 2219.7292
$ws.Range("I134").Value = 2282.8865
$ws.Range("K134").Value = 6848.6595
$ws.Range("M134").Value = -4313.6595

$ws = $wb.Worksheets.Item("CRP")
# Row 15: On the Move
$ws.Range("H15").Value = 5648.625
$ws.Range("I15").Value = 2478
$ws.Range("J15").Value = 10933
$ws.Range("K15").Value = 2478
$ws.Range("L15").Value = 10933
$ws.Range("M15").Value = -2308
$ws.Range("N15").Value = -11273

# Row 22: Driving Up the Wall
$ws.Range("H22").Value = 3880.7144
$ws.Range("I22").Value = 3721.6667
$ws.Range("J22").Value = 4000
$ws.Range("K22").Value = 3721.6667
$ws.Range("L22").Value = 4000
$ws.Range("M22").Value = -3371.6667
$ws.Range("N22").Value = -4700

# Row 31: Wall Not Found
$ws.Range("H31").Value = 6769.8125
$ws.Range("I31").Value = 5700
$ws.Range("J31").Value = 8552.833000000001
$ws.Range("K31").Value = 5700
$ws.Range("L31").Value = 8552.833000000001
$ws.Range("M31").Value = -5405
$ws.Range("N31").Value = -9142.833000000001

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 6769.8125
$ws.Range("I34").Value = 5700
$ws.Range("J34").Value = 8552.833000000001
$ws.Range("K34").Value = 5700
$ws.Range("L34").Value = 8552.833000000001
$ws.Range("M34").Value = -5498
$ws.Range("N34").Value = -8956.833000000001

# Row 35: Storm of Swords
$ws.Range("H35").Value = 3475.875
$ws.Range("I35").Value = 1967.8334
$ws.Range("K35").Value = 1967.8334
$ws.Range("M35").Value = -1673.8334

# Row 99: O Pine
$ws.Range("H99").Value = 3907.1428
$ws.Range("I99").Value = 4141.6665
$ws.Range("K99").Value = 4141.6665
$ws.Range("M99").Value = -2643.6665

# Row 126: A Better Conductor
$ws.Range("H126").Value = 3907.1428
$ws.Range("I126").Value = 4141.6665
$ws.Range("K126").Value = 12424.9995
$ws.Range("M126").Value = -9954.999500000002

# Row 139: Weaving a Path
$ws.Range("H139").Value = 144722.5
$ws.Range("I139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("M139").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 34: Fever Pitch
$ws.Range("H34").Value = 449.36365
$ws.Range("I34").Value = 183.33333
$ws.Range("J34").Value = 549.125
$ws.Range("K34").Value = 549.99999
$ws.Range("L34").Value = 1647.375
$ws.Range("M34").Value = -465.99999
$ws.Range("N34").Value = -1815.375

# Row 37: I Love Lamprey
$ws.Range("H37").Value = 58227.58
$ws.Range("J37").Value = 58227.58
$ws.Range("L37").Value = 174682.74
$ws.Range("N37").Value = -174906.74

# Row 38: Pretty as a Picture
$ws.Range("H38").Value = 713.75
$ws.Range("I38").Value = 26.666666
$ws.Range("K38").Value = 79.99999800000001
$ws.Range("M38").Value = 267.000002

# Row 39: Bloody Good Tart, This
$ws.Range("H39").Value = 1333.3334

# Row 55: Pagan Pastries
$ws.Range("H55").Value = 877.3077
$ws.Range("J55").Value = 999.8182
$ws.Range("L55").Value = 2999.4546
$ws.Range("N55").Value = -3353.4546

$ws = $wb.Worksheets.Item("GSM")
# Row 99: Needle in a Hingan Stack
$ws.Range("H99").Value = 12335
$ws.Range("I99").Value = 7202
$ws.Range("K99").Value = 7202
$ws.Range("M99").Value = -4956

# Row 132: On Board for Lar
$ws.Range("H132").Value = 2736.487
$ws.Range("I132").Value = 2391.2666
$ws.Range("J132").Value = 3887.2222
$ws.Range("K132").Value = 7173.7998
$ws.Range("L132").Value = 11661.6666
$ws.Range("M132").Value = -4643.7998
$ws.Range("N132").Value = -16721.6666

$ws = $wb.Worksheets.Item("LTW")
# Row 16: Saddle Sore
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()

# Row 61: Spelling Me Softly
$ws.Range("H61").Value = 59402.277
$ws.Range("I61").Value = 70249.39999999999
$ws.Range("J61").Value = 5166.6665
$ws.Range("K61").Value = 70249.39999999999
$ws.Range("L61").Value = 5166.6665
$ws.Range("M61").Value = -70047.39999999999
$ws.Range("N61").Value = -5570.6665

# Row 113: Peace in Rest
$ws.Range("H113").Value = 59402.277
$ws.Range("I113").Value = 70249.39999999999
$ws.Range("J113").Value = 5166.6665
$ws.Range("K113").Value = 70249.39999999999
$ws.Range("L113").Value = 5166.6665
$ws.Range("M113").Value = -68079.39999999999
$ws.Range("N113").Value = -9506.666499999999

# Row 122: Hell on Leather
$ws.Range("H122").Value = 4000
$ws.Range("I122").Value = 3000
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = -19900

$ws = $wb.Worksheets.Item("WVR")
# Row 81: Where the Dragonflies, the Net Catches
$ws.Range("H81").Value = 1913.5714
$ws.Range("I81").Value = 1913.5714
$ws.Range("K81").Value = 3827.1428
$ws.Range("M81").Value = -2766.1428

# Row 84: To Kill a Dragon on Nameday (L)
$ws.Range("H84").Value = 1913.5714
$ws.Range("I84").Value = 1913.5714
$ws.Range("K84").Value = 19135.714
$ws.Range("M84").Value = -13831.714

# Row 126: A Polished Purchase
$ws.Range("H126").Value = 2848.7
$ws.Range("I126").Value = 3214.2856
$ws.Range("K126").Value = 9642.856800000001
$ws.Range("M126").Value = -7172.856800000001
